# Added USDA thermal curve data: the "plate" column (C) held the plate
# number as a bare numeric value (1 or 2). Re-express it as a descriptive
# text label ("plate1" / "plate2") so it reads the same as the other
# category columns (well, etc.) in this worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-25: plate 2 -> "plate2"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = "plate2"
}

# Rows 26-97: plate 1 -> "plate1"
for ($r = 26; $r -le 97; $r++) {
    $ws.Cells.Item($r, 3).Value = "plate1"
}

# Restore the author's last on-screen selection.
$ws.Range("C89").Select()

Write-Host "plate column relabeled as text"
